$d = $word.ActiveDocument

# 1. Replace the title text "2.2 - Debate I" with "Placeholder - Check Back Later".
#    Assign directly to the paragraph's Range.Text (rather than Find/Replace) so
#    the run keeps its normal xml:space="preserve" text serialization.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Title" -and $p.Range.Text -like "*2.2 - Debate I*") {
        $p.Range.Text = "Placeholder - Check Back Later"
        break
    }
}

# 2. Remove the trailing " " and ":::" runs that follow the
#    "do not need to be looked at..." sentence in the Additional Resources cell.
#    Locate the " :::" span with Find (no replacement text) and delete that
#    range outright so the preceding run is left untouched.
$r = $d.Content
$found = $r.Find.Execute(" :::", $false, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
if ($found) {
    $r.Delete()
}
